$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume (columns D/E) hold free-form text ("29.509.84", "1.000",
# "  -0.82%  ", ...). Excel's Value setter auto-coerces numeric-looking
# strings, so each D/E cell we touch is switched to Text format ("@")
# first to preserve the literal text, exactly like the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.509.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.850.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.26'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6372'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07569'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2994'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.21'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07687'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.863.09'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.029'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6873'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '84.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009806'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.110.13'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.253'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.12%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.550.48'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.66%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.47'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.53'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.628'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.54%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.30'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1393'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.456'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.77'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.489'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05892'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.29%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.276'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.128'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.58%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.062'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.907'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.172'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7188'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.598'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.803'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.237.52'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01776'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9146'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.123'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.026.28'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.80'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.75%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.384'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +10.38%  '
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4040'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.155'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.67%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.703'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.72%  '
